$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 7)
$ws.Range("A7").Value = "https://www.smarttask.io/"
$ws.Range("B7").Value = "smarttask.io/"

# Turn A7 into a hyperlink (this also applies the built-in "Hyperlink" style)
$ws.Hyperlinks.Add($ws.Range("A7"), "https://www.smarttask.io/") | Out-Null

# Column A is now the widest (URLs), so Excel "best fit" widens it
$ws.Columns("A:A").ColumnWidth = 32.8

# Select B7 and mark the sheet as the active tab (matches final sheetView state)
$ws.Range("B7").Select()
